$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 1265.9166
$ws.Range("J19").Value = 1531.8334
$ws.Range("L19").Value = 1531.8334
$ws.Range("N19").Value = -1881.8334
$ws.Range("H69").Value = 9988.571
$ws.Range("J69").Value = 10014.588
$ws.Range("L69").Value = 30043.764
$ws.Range("N69").Value = -31791.764
$ws.Range("H70").Value = 83341660
$ws.Range("I70").Value = 0
$ws.Range("J70").Value = 83341660
$ws.Range("K70").Value = 0
$ws.Range("L70").Value = 250024980
$ws.Range("M70").Value = $null
$ws.Range("N70").Value = -250025520
$ws.Range("H72").Value = 9988.571
$ws.Range("J72").Value = 10014.588
$ws.Range("L72").Value = 90131.292
$ws.Range("N72").Value = -98867.292
$ws.Range("H73").Value = 83341660
$ws.Range("I73").Value = 0
$ws.Range("J73").Value = 83341660
$ws.Range("K73").Value = 0
$ws.Range("L73").Value = 250024980
$ws.Range("M73").Value = $null
$ws.Range("N73").Value = -250026852
$ws.Range("H76").Value = 5685.077
$ws.Range("I76").Value = 5667.6665
$ws.Range("J76").Value = 5894
$ws.Range("K76").Value = 5667.6665
$ws.Range("L76").Value = 5894
$ws.Range("M76").Value = -5352.6665
$ws.Range("N76").Value = -6524
$ws.Range("H79").Value = 5685.077
$ws.Range("I79").Value = 5667.6665
$ws.Range("J79").Value = 5894
$ws.Range("K79").Value = 5667.6665
$ws.Range("L79").Value = 5894
$ws.Range("M79").Value = -4575.6665
$ws.Range("N79").Value = -8078
$ws.Range("H82").Value = 2933
$ws.Range("J82").Value = 4000
$ws.Range("L82").Value = 12000
$ws.Range("N82").Value = -12812
$ws.Range("H85").Value = 2933
$ws.Range("J85").Value = 4000
$ws.Range("L85").Value = 12000
$ws.Range("N85").Value = -14808
$ws.Range("H98").Value = 1603.35
$ws.Range("I98").Value = 1603.35
$ws.Range("K98").Value = 1603.35
$ws.Range("M98").Value = -105.3499999999999
$ws.Range("H107").Value = 2774.7273
$ws.Range("I107").Value = 3985.7144
$ws.Range("K107").Value = 3985.7144
$ws.Range("M107").Value = -2065.7144
$ws.Range("H112").Value = 1079031.4
$ws.Range("J112").Value = 1186742.1
$ws.Range("L112").Value = 3560226.3
$ws.Range("N112").Value = -3562442.3
$ws.Range("H122").Value = 1603.35
$ws.Range("I122").Value = 1603.35
$ws.Range("K122").Value = 4810.049999999999
$ws.Range("M122").Value = -2360.049999999999
$ws.Range("H132").Value = 3703.825
$ws.Range("I132").Value = 3744.6
$ws.Range("K132").Value = 11233.8
$ws.Range("M132").Value = -8703.799999999999
$ws.Range("H135").Value = 928.4167
$ws.Range("I135").Value = 928.4167
$ws.Range("K135").Value = 8355.7503
$ws.Range("M135").Value = -5820.7503
$ws.Range("H137").Value = 1761.6274
$ws.Range("I137").Value = 1975.8889
$ws.Range("J137").Value = 1644.7576
$ws.Range("K137").Value = 5927.6667
$ws.Range("L137").Value = 4934.2728
$ws.Range("M137").Value = -3377.6667
$ws.Range("N137").Value = -10034.2728
$ws.Range("H138").Value = 3139.6453
$ws.Range("I138").Value = 2403.697
$ws.Range("J138").Value = 3544.4167
$ws.Range("K138").Value = 7211.091
$ws.Range("L138").Value = 10633.2501
$ws.Range("M138").Value = -2071.091
$ws.Range("N138").Value = -20913.2501
$ws.Range("H140").Value = 255799
$ws.Range("J140").Value = 277250
$ws.Range("L140").Value = 277250
$ws.Range("N140").Value = -287610
$ws.Range("H141").Value = 4910.636
$ws.Range("I141").Value = 5001.7
$ws.Range("K141").Value = 15005.1
$ws.Range("M141").Value = -9825.099999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2506542.2
$ws.Range("I32").Value = 2814349.8
$ws.Range("J32").Value = 16101.091
$ws.Range("K32").Value = 2814349.8
$ws.Range("L32").Value = 16101.091
$ws.Range("M32").Value = -2814062.8
$ws.Range("N32").Value = -16675.091
$ws.Range("H45").Value = 3400.261
$ws.Range("I45").Value = 3094.7058
$ws.Range("J45").Value = 4266
$ws.Range("K45").Value = 3094.7058
$ws.Range("L45").Value = 4266
$ws.Range("M45").Value = -2717.7058
$ws.Range("N45").Value = -5020
$ws.Range("H61").Value = 3047.6042
$ws.Range("I61").Value = 3039.75
$ws.Range("J61").Value = 3071.1667
$ws.Range("K61").Value = 3039.75
$ws.Range("L61").Value = 3071.1667
$ws.Range("M61").Value = -2827.75
$ws.Range("N61").Value = -3495.1667
$ws.Range("H74").Value = 2802.9092
$ws.Range("I74").Value = 2532.875
$ws.Range("K74").Value = 2532.875
$ws.Range("M74").Value = -1658.875
$ws.Range("H77").Value = 2802.9092
$ws.Range("I77").Value = 2532.875
$ws.Range("K77").Value = 12664.375
$ws.Range("M77").Value = -8296.375
$ws.Range("H132").Value = 241256.73
$ws.Range("I132").Value = 437155.25
$ws.Range("K132").Value = 1311465.75
$ws.Range("M132").Value = -1308935.75
$ws.Range("H136").Value = 3047.6042
$ws.Range("I136").Value = 3039.75
$ws.Range("J136").Value = 3071.1667
$ws.Range("K136").Value = 9119.25
$ws.Range("L136").Value = 9213.500100000001
$ws.Range("M136").Value = -6569.25
$ws.Range("N136").Value = -14313.5001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H117").Value = 123000
$ws.Range("J117").Value = 123000
$ws.Range("L117").Value = 123000
$ws.Range("N117").Value = -132178
$ws.Range("H134").Value = 9414870
$ws.Range("I134").Value = 1932000.6
$ws.Range("K134").Value = 5796001.800000001
$ws.Range("M134").Value = -5793466.800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 4808.6206
$ws.Range("I31").Value = 3500.3809
$ws.Range("K31").Value = 3500.3809
$ws.Range("M31").Value = -3205.3809
$ws.Range("H34").Value = 4808.6206
$ws.Range("I34").Value = 3500.3809
$ws.Range("K34").Value = 3500.3809
$ws.Range("M34").Value = -3298.3809
$ws.Range("H58").Value = 3537.4666
$ws.Range("I58").Value = 2985.2632
$ws.Range("K58").Value = 2985.2632
$ws.Range("M58").Value = -2782.2632
$ws.Range("H132").Value = 1274.5
$ws.Range("I132").Value = 1274.5
$ws.Range("K132").Value = 3823.5
$ws.Range("M132").Value = -1293.5
$ws.Range("H134").Value = 2412.3333
$ws.Range("I134").Value = 1900
$ws.Range("K134").Value = 5700
$ws.Range("M134").Value = -3165
$ws.Range("H136").Value = 3537.4666
$ws.Range("I136").Value = 2985.2632
$ws.Range("K136").Value = 8955.7896
$ws.Range("M136").Value = -6405.7896

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 121049384
$ws.Range("I4").Value = 92432104
$ws.Range("J4").Value = 199746910
$ws.Range("K4").Value = 277296312
$ws.Range("L4").Value = 599240730
$ws.Range("M4").Value = -277296200
$ws.Range("N4").Value = -599240954
$ws.Range("H68").Value = 784
$ws.Range("I68").Value = 851
$ws.Range("K68").Value = 2553
$ws.Range("M68").Value = -1742
$ws.Range("H71").Value = 784
$ws.Range("I71").Value = 851
$ws.Range("K71").Value = 7659
$ws.Range("M71").Value = -3603
$ws.Range("I140").Value = 11767561
$ws.Range("K140").Value = 35302683
$ws.Range("M140").Value = -35297503

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3604.1052
$ws.Range("I132").Value = 2958.7693
$ws.Range("J132").Value = 5002.3335
$ws.Range("K132").Value = 8876.3079
$ws.Range("L132").Value = 15007.0005
$ws.Range("M132").Value = -6346.3079
$ws.Range("N132").Value = -20067.0005

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 2288.5
$ws.Range("I7").Value = 2288.5
$ws.Range("J7").Value = 0
$ws.Range("K7").Value = 2288.5
$ws.Range("L7").Value = 0
$ws.Range("M7").Value = -2176.5
$ws.Range("N7").Value = $null
$ws.Range("H122").Value = 4313908.5
$ws.Range("I122").Value = 5750210.5
$ws.Range("K122").Value = 17250631.5
$ws.Range("M122").Value = -17248181.5
$ws.Range("H126").Value = 2288.5
$ws.Range("I126").Value = 2288.5
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 6865.5
$ws.Range("L126").Value = 0
$ws.Range("M126").Value = -4395.5
$ws.Range("N126").Value = $null
$ws.Range("H132").Value = 125588.83
$ws.Range("I132").Value = 315473.66
$ws.Range("K132").Value = 946420.98
$ws.Range("M132").Value = -943890.98
$ws.Range("H136").Value = 3671.1775
$ws.Range("I136").Value = 3342.0186
$ws.Range("K136").Value = 10026.0558
$ws.Range("M136").Value = -7476.0558
$ws.Range("H140").Value = 86199.2
$ws.Range("I140").Value = 79500
$ws.Range("K140").Value = 79500
$ws.Range("M140").Value = -74320

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 4068.3333
$ws.Range("I6").Value = 1205
$ws.Range("J6").Value = 5500
$ws.Range("K6").Value = 1205
$ws.Range("L6").Value = 5500
$ws.Range("M6").Value = -1090
$ws.Range("N6").Value = -5730
$ws.Range("H126").Value = 8832.462
$ws.Range("I126").Value = 8534.111000000001
$ws.Range("K126").Value = 25602.333
$ws.Range("M126").Value = -23132.333
$ws.Range("H132").Value = 28016.775
$ws.Range("I132").Value = 39621.965
$ws.Range("K132").Value = 118865.895
$ws.Range("M132").Value = -116335.895
$ws.Range("H136").Value = 26772.822
$ws.Range("I136").Value = 2175.5588
$ws.Range("J136").Value = 102800.73
$ws.Range("K136").Value = 6526.676399999999
$ws.Range("L136").Value = 308402.19
$ws.Range("M136").Value = -3976.676399999999
$ws.Range("N136").Value = -313502.19
$ws.Range("H137").Value = 74157.60000000001
$ws.Range("J137").Value = 74157.60000000001
$ws.Range("L137").Value = 74157.60000000001
$ws.Range("N137").Value = -84357.60000000001

